# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt" /
# "Arándano (blue)" right above the current row 70, pushing the existing
# rows 70-76 down to 71-77.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 70 (shifts old rows 70-76 -> 71-77)
$ws.Rows.Item(70).Insert()

# Populate the new row 70 with the new weekly record
$ws.Cells.Item(70, 1).Value = 4
$ws.Cells.Item(70, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(70, 3).Value = "Los Lagos"
$ws.Cells.Item(70, 4).Value = 45265
$ws.Cells.Item(70, 5).Value = 10
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100101
$ws.Cells.Item(70, 8).Value = "Berries"
$ws.Cells.Item(70, 9).Value = 100101001
$ws.Cells.Item(70, 10).Value = "Arándano (blue)"
$ws.Cells.Item(70, 11).Value = "Sin especificar"
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 300
$ws.Cells.Item(70, 14).Value = 6000
$ws.Cells.Item(70, 15).Value = 6000
$ws.Cells.Item(70, 16).Value = 6000
$ws.Cells.Item(70, 17).Value = "$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(70, 18).Value = "Región del Maule"
$ws.Cells.Item(70, 19).Value = 4000
$ws.Cells.Item(70, 20).Value = 1.5
